$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logBook")

# New row 39 data - entries for 2nd July continued until 7:30pm
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = 44744
$ws.Cells.Item(39, 3).Value = 0.66666666666666663
$ws.Cells.Item(39, 4).Value = 0.8125
$ws.Cells.Item(39, 5).Formula = "=D39-C39"
$ws.Cells.Item(39, 6).Value = "Code"
$ws.Cells.Item(39, 7).Value = "1. Formatted and uploaded FCN_resnet50_baseline nb`n2. Unet_starter nb completed`n3. Unet_baseline nb completed"

# Copy formatting from row 38 to new row 39 (styles for each cell A-G)
$ws.Range("A38:G38").Copy() | Out-Null
$ws.Range("A39:G39").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-apply the values/formula after paste-special formats (paste formats only shouldn't touch values,
# but ensure everything is correct regardless)
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = 44744
$ws.Cells.Item(39, 3).Value = 0.66666666666666663
$ws.Cells.Item(39, 4).Value = 0.8125
$ws.Cells.Item(39, 5).Formula = "=D39-C39"
$ws.Cells.Item(39, 6).Value = "Code"
$ws.Cells.Item(39, 7).Value = "1. Formatted and uploaded FCN_resnet50_baseline nb`n2. Unet_starter nb completed`n3. Unet_baseline nb completed"

$ws.Rows.Item(39).RowHeight = 45

# Update selection to E50 (as if user clicked the total cell after adding the entry)
$ws.Range("E50").Select() | Out-Null

$wb.Save()
